$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1760.6
$ws.Range("I31").Value = 767.6667
$ws.Range("J31").Value = 3250
$ws.Range("K31").Value = 2303.0001
$ws.Range("L31").Value = 9750
$ws.Range("M31").Value = -2073.0001
$ws.Range("N31").Value = -10210

$ws.Range("H33").Value = 223
$ws.Range("I33").Value = 96.30768999999999
$ws.Range("K33").Value = 96.30768999999999
$ws.Range("M33").Value = 132.69231

$ws.Range("H86").Value = 4187.8125
$ws.Range("I86").Value = 3400.3845
$ws.Range("J86").Value = 7600
$ws.Range("K86").Value = 3400.3845
$ws.Range("L86").Value = 7600
$ws.Range("M86").Value = -2277.3845
$ws.Range("N86").Value = -9846

$ws.Range("H89").Value = 4187.8125
$ws.Range("I89").Value = 3400.3845
$ws.Range("J89").Value = 7600
$ws.Range("K89").Value = 17001.9225
$ws.Range("L89").Value = 38000
$ws.Range("M89").Value = -11385.9225
$ws.Range("N89").Value = -49232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8100.65
$ws.Range("I32").Value = 7423.2163
$ws.Range("J32").Value = 30004.334
$ws.Range("K32").Value = 7423.2163
$ws.Range("L32").Value = 30004.334
$ws.Range("M32").Value = -7136.2163
$ws.Range("N32").Value = -30578.334

$ws.Range("H74").Value = 1103.5135
$ws.Range("I74").Value = 1233.44
$ws.Range("J74").Value = 832.8333
$ws.Range("K74").Value = 1233.44
$ws.Range("L74").Value = 832.8333
$ws.Range("M74").Value = -359.4400000000001
$ws.Range("N74").Value = -2580.8333

$ws.Range("H77").Value = 1103.5135
$ws.Range("I77").Value = 1233.44
$ws.Range("J77").Value = 832.8333
$ws.Range("K77").Value = 6167.200000000001
$ws.Range("L77").Value = 4164.1665
$ws.Range("M77").Value = -1799.200000000001
$ws.Range("N77").Value = -12900.1665

$ws.Range("H88").Value = 2990.5
$ws.Range("I88").Value = 2246.5
$ws.Range("J88").Value = 3734.5
$ws.Range("K88").Value = 2246.5
$ws.Range("L88").Value = 3734.5
$ws.Range("M88").Value = -1840.5
$ws.Range("N88").Value = -4546.5

$ws.Range("H91").Value = 2990.5
$ws.Range("I91").Value = 2246.5
$ws.Range("J91").Value = 3734.5
$ws.Range("K91").Value = 2246.5
$ws.Range("L91").Value = 3734.5
$ws.Range("M91").Value = -842.5
$ws.Range("N91").Value = -6542.5

$ws.Range("H125").Value = 28999.111
$ws.Range("J125").Value = 28999.111
$ws.Range("L125").Value = 28999.111
$ws.Range("N125").Value = -38839.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 28000
$ws.Range("J62").Value = 28000
$ws.Range("L62").Value = 28000
$ws.Range("N62").Value = -29372

$ws.Range("H65").Value = 28000
$ws.Range("J65").Value = 28000
$ws.Range("L65").Value = 84000
$ws.Range("N65").Value = -90864

$ws.Range("H94").Value = 1269
$ws.Range("I94").Value = 1273.1482
$ws.Range("J94").Value = 1255
$ws.Range("K94").Value = 1273.1482
$ws.Range("L94").Value = 1255
$ws.Range("M94").Value = -822.1482000000001
$ws.Range("N94").Value = -2157

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27699.896
$ws.Range("I31").Value = 1965.84
$ws.Range("J31").Value = 73653.57000000001
$ws.Range("K31").Value = 1965.84
$ws.Range("L31").Value = 73653.57000000001
$ws.Range("M31").Value = -1670.84
$ws.Range("N31").Value = -74243.57000000001

$ws.Range("H34").Value = 27699.896
$ws.Range("I34").Value = 1965.84
$ws.Range("J34").Value = 73653.57000000001
$ws.Range("K34").Value = 1965.84
$ws.Range("L34").Value = 73653.57000000001
$ws.Range("M34").Value = -1763.84
$ws.Range("N34").Value = -74057.57000000001

$ws.Range("H99").Value = 1935.1
$ws.Range("I99").Value = 1934.6207
$ws.Range("J99").Value = 1936.3636
$ws.Range("K99").Value = 1934.6207
$ws.Range("L99").Value = 1936.3636
$ws.Range("M99").Value = -436.6206999999999
$ws.Range("N99").Value = -4932.3636

$ws.Range("H107").Value = 393.68182
$ws.Range("I107").Value = 364.1
$ws.Range("J107").Value = 418.33334
$ws.Range("K107").Value = 364.1
$ws.Range("L107").Value = 418.33334
$ws.Range("M107").Value = 1555.9
$ws.Range("N107").Value = -4258.33334

$ws.Range("H126").Value = 1935.1
$ws.Range("I126").Value = 1934.6207
$ws.Range("J126").Value = 1936.3636
$ws.Range("K126").Value = 5803.8621
$ws.Range("L126").Value = 5809.0908
$ws.Range("M126").Value = -3333.8621
$ws.Range("N126").Value = -10749.0908

$ws.Range("H132").Value = 1926.0588
$ws.Range("I132").Value = 1487.875
$ws.Range("K132").Value = 4463.625
$ws.Range("M132").Value = -1933.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 713.2308
$ws.Range("I113").Value = 864
$ws.Range("J113").Value = 633.41174
$ws.Range("K113").Value = 2592
$ws.Range("L113").Value = 1900.23522
$ws.Range("M113").Value = -422
$ws.Range("N113").Value = -6240.23522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2769.5
$ws.Range("I126").Value = 3073.818
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 9221.454000000002
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -6751.454000000002
$ws.Range("N126").Value = -11240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1872.76
$ws.Range("I61").Value = 1871.0555
$ws.Range("J61").Value = 1877.1428
$ws.Range("K61").Value = 1871.0555
$ws.Range("L61").Value = 1877.1428
$ws.Range("M61").Value = -1669.0555
$ws.Range("N61").Value = -2281.1428

$ws.Range("H106").Value = 17500
$ws.Range("J106").Value = 17500
$ws.Range("L106").Value = 17500
$ws.Range("N106").Value = -20024

$ws.Range("H113").Value = 1872.76
$ws.Range("I113").Value = 1871.0555
$ws.Range("J113").Value = 1877.1428
$ws.Range("K113").Value = 1871.0555
$ws.Range("L113").Value = 1877.1428
$ws.Range("M113").Value = 298.9445000000001
$ws.Range("N113").Value = -6217.1428

$ws.Range("H132").Value = 3205.6287
$ws.Range("I132").Value = 2105.2104
$ws.Range("K132").Value = 6315.6312
$ws.Range("M132").Value = -3785.6312

$ws.Range("H136").Value = 3489.0754
$ws.Range("I136").Value = 1747.9
$ws.Range("J136").Value = 8846.538
$ws.Range("K136").Value = 5243.700000000001
$ws.Range("L136").Value = 26539.614
$ws.Range("M136").Value = -2693.700000000001
$ws.Range("N136").Value = -31639.614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14670.8
$ws.Range("J41").Value = 14670.8
$ws.Range("L41").Value = 14670.8
$ws.Range("N41").Value = -15450.8

$ws.Range("H101").Value = 15574.333
$ws.Range("J101").Value = 15574.333
$ws.Range("L101").Value = 15574.333
$ws.Range("N101").Value = -22064.333

$ws.Range("H132").Value = 824.18054
$ws.Range("I132").Value = 652.9820999999999
$ws.Range("J132").Value = 1423.375
$ws.Range("K132").Value = 1958.9463
$ws.Range("L132").Value = 4270.125
$ws.Range("M132").Value = 571.0537000000002
$ws.Range("N132").Value = -9330.125
